# Update row 2 of Sheet1 ("利润表/688658.xlsx") from the FY2019 annual
# report figures to the 2020-06-30 half-year report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: 001 -> 002 (must remain a text value, not become numeric 2)
$ws.Range("J2").Value = "'002"

# REPORT_DATE: 2019-12-31 -> 2020-06-30
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures that get new values
$ws.Range("O2").Value = 132251328.35
$ws.Range("P2").Value = 1681584033.13
$ws.Range("Q2").Value = 1541694663.03
$ws.Range("S2").Value = 621242077.4
$ws.Range("T2").Value = 621242077.4
$ws.Range("V2").Value = 741792240.2
$ws.Range("W2").Value = 85646338.57
$ws.Range("X2").Value = 20930930.43
$ws.Range("Y2").Value = 145884910.55
$ws.Range("Z2").Value = 147439790.61
$ws.Range("AA2").Value = 15142154.37
$ws.Range("AG2").Value = 18290519.3
$ws.Range("AS2").Value = 120816425.8

# Cells that previously held computed ratios but are now blanked out
# (kept present in the sheet as empty cells, matching the source data feed).
$blankCells = @("R2", "U2", "AP2", "AQ2", "AR2", "AT2")
foreach ($addr in $blankCells) {
    $rng = $ws.Range($addr)
    $rng.Value = ""
    $rng.Style = "Normal"
}
